$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 874.875
$ws.Range("I6").Value = 857
$ws.Range("K6").Value = 2571
$ws.Range("M6").Value = -2459
$ws.Range("H32").Value = 22728378
$ws.Range("I32").Value = 108
$ws.Range("K32").Value = 108
$ws.Range("M32").Value = 218
$ws.Range("H33").Value = 5130442
$ws.Range("I33").Value = 6839618.5
$ws.Range("J33").Value = 2911.75
$ws.Range("K33").Value = 6839618.5
$ws.Range("L33").Value = 2911.75
$ws.Range("M33").Value = -6839389.5
$ws.Range("N33").Value = -3369.75
$ws.Range("H76").Value = 36365.832
$ws.Range("I76").Value = 38641.965
$ws.Range("J76").Value = 4500
$ws.Range("K76").Value = 38641.965
$ws.Range("L76").Value = 4500
$ws.Range("M76").Value = -38326.965
$ws.Range("N76").Value = -5130
$ws.Range("H79").Value = 36365.832
$ws.Range("I79").Value = 38641.965
$ws.Range("J79").Value = 4500
$ws.Range("K79").Value = 38641.965
$ws.Range("L79").Value = 4500
$ws.Range("M79").Value = -37549.965
$ws.Range("N79").Value = -6684
$ws.Range("H107").Value = 2937.6
$ws.Range("I107").Value = 2937.6
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2937.6
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1017.6
$ws.Range("N107").ClearContents()
$ws.Range("H132").Value = 10405731
$ws.Range("I132").Value = 13159393
$ws.Range("K132").Value = 39478179
$ws.Range("M132").Value = -39475649
$ws.Range("H135").Value = 710.7143
$ws.Range("I135").Value = 726.9091
$ws.Range("K135").Value = 6542.1819
$ws.Range("M135").Value = -4007.1819
$ws.Range("H137").Value = 1593457.8
$ws.Range("I137").Value = 5151.7427
$ws.Range("K137").Value = 15455.2281
$ws.Range("M137").Value = -12905.2281

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 2886
$ws.Range("I22").Value = 1863.2
$ws.Range("J22").Value = 8000
$ws.Range("K22").Value = 1863.2
$ws.Range("L22").Value = 8000
$ws.Range("M22").Value = -1564.2
$ws.Range("N22").Value = -8598
$ws.Range("H29").Value = 33336224
$ws.Range("I29").Value = 33336224
$ws.Range("K29").Value = 33336224
$ws.Range("M29").Value = -33335916
$ws.Range("H61").Value = 1429164.9
$ws.Range("I61").Value = 48256.477
$ws.Range("J61").Value = 3100790.8
$ws.Range("K61").Value = 48256.477
$ws.Range("L61").Value = 3100790.8
$ws.Range("M61").Value = -48044.477
$ws.Range("N61").Value = -3101214.8
$ws.Range("H74").Value = 394485.16
$ws.Range("I74").Value = 3432.7273
$ws.Range("K74").Value = 3432.7273
$ws.Range("M74").Value = -2558.7273
$ws.Range("H77").Value = 394485.16
$ws.Range("I77").Value = 3432.7273
$ws.Range("K77").Value = 17163.6365
$ws.Range("M77").Value = -12795.6365
$ws.Range("H102").Value = 2092
$ws.Range("J102").Value = 1758.6
$ws.Range("L102").Value = 1758.6
$ws.Range("N102").Value = -5002.6
$ws.Range("H136").Value = 1429164.9
$ws.Range("I136").Value = 48256.477
$ws.Range("J136").Value = 3100790.8
$ws.Range("K136").Value = 144769.431
$ws.Range("L136").Value = 9302372.399999999
$ws.Range("M136").Value = -142219.431
$ws.Range("N136").Value = -9307472.399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 12717.72
$ws.Range("I107").Value = 14722.263
$ws.Range("K107").Value = 14722.263
$ws.Range("M107").Value = -12802.263

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4382.921
$ws.Range("I31").Value = 4420.5
$ws.Range("K31").Value = 4420.5
$ws.Range("M31").Value = -4125.5
$ws.Range("H34").Value = 4382.921
$ws.Range("I34").Value = 4420.5
$ws.Range("K34").Value = 4420.5
$ws.Range("M34").Value = -4218.5
$ws.Range("H58").Value = 1839.1
$ws.Range("I58").Value = 1578.2858
$ws.Range("K58").Value = 1578.2858
$ws.Range("M58").Value = -1375.2858
$ws.Range("H94").Value = 983.3158
$ws.Range("I94").Value = 931.4286
$ws.Range("K94").Value = 931.4286
$ws.Range("M94").Value = -480.4286
$ws.Range("H134").Value = 2938.1
$ws.Range("I134").Value = 2422.75
$ws.Range("K134").Value = 7268.25
$ws.Range("M134").Value = -4733.25
$ws.Range("H136").Value = 1839.1
$ws.Range("I136").Value = 1578.2858
$ws.Range("K136").Value = 4734.857400000001
$ws.Range("M136").Value = -2184.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 323.4
$ws.Range("I11").Value = 280.5
$ws.Range("J11").Value = 399.66666
$ws.Range("K11").Value = 841.5
$ws.Range("L11").Value = 1198.99998
$ws.Range("M11").Value = -701.5
$ws.Range("N11").Value = -1478.99998
$ws.Range("H125").Value = 3499
$ws.Range("I125").Value = 5000
$ws.Range("J125").Value = 1998
$ws.Range("K125").Value = 15000
$ws.Range("L125").Value = 5994
$ws.Range("M125").Value = -10080
$ws.Range("N125").Value = -15834
$ws.Range("H131").Value = 4042754.2
$ws.Range("I131").Value = 5682896
$ws.Range("K131").Value = 17048688
$ws.Range("M131").Value = -17043648

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 29414106
$ws.Range("I102").Value = 38463344
$ws.Range("K102").Value = 38463344
$ws.Range("M102").Value = -38461722
$ws.Range("H105").Value = 51223.75
$ws.Range("J105").Value = 51223.75
$ws.Range("L105").Value = 51223.75
$ws.Range("N105").Value = -58211.75
$ws.Range("H132").Value = 6893968.5
$ws.Range("I132").Value = 2838.682
$ws.Range("J132").Value = 15316461
$ws.Range("K132").Value = 8516.045999999998
$ws.Range("L132").Value = 45949383
$ws.Range("M132").Value = -5986.045999999998
$ws.Range("N132").Value = -45954443

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 34411.766
$ws.Range("I20").Value = 44583.332
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 44583.332
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = -44357.332
$ws.Range("N20").Value = -10452
$ws.Range("H22").Value = 2627.5356
$ws.Range("I22").Value = 1058.4
$ws.Range("J22").Value = 2968.652
$ws.Range("K22").Value = 1058.4
$ws.Range("L22").Value = 2968.652
$ws.Range("M22").Value = -763.4000000000001
$ws.Range("N22").Value = -3558.652
$ws.Range("H27").Value = 2627.5356
$ws.Range("I27").Value = 1058.4
$ws.Range("J27").Value = 2968.652
$ws.Range("K27").Value = 1058.4
$ws.Range("L27").Value = 2968.652
$ws.Range("M27").Value = -951.4000000000001
$ws.Range("N27").Value = -3182.652
$ws.Range("H61").Value = 2759.9487
$ws.Range("I61").Value = 1719.8788
$ws.Range("K61").Value = 1719.8788
$ws.Range("M61").Value = -1517.8788
$ws.Range("H93").Value = 1393.6364
$ws.Range("I93").Value = 1413
$ws.Range("J93").Value = 1200
$ws.Range("K93").Value = 1413
$ws.Range("L93").Value = 1200
$ws.Range("M93").Value = -165
$ws.Range("N93").Value = -3696
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("H106").Value = 16547.6
$ws.Range("J106").Value = 16547.6
$ws.Range("L106").Value = 16547.6
$ws.Range("N106").Value = -19071.6
$ws.Range("H113").Value = 2759.9487
$ws.Range("I113").Value = 1719.8788
$ws.Range("K113").Value = 1719.8788
$ws.Range("M113").Value = 450.1212
$ws.Range("H132").Value = 3132.457
$ws.Range("I132").Value = 2556.4583
$ws.Range("K132").Value = 7669.374899999999
$ws.Range("M132").Value = -5139.374899999999
$ws.Range("H136").Value = 33385.4
$ws.Range("I136").Value = 40245.18
$ws.Range("J136").Value = 5946.2856
$ws.Range("K136").Value = 120735.54
$ws.Range("L136").Value = 17838.8568
$ws.Range("M136").Value = -118185.54
$ws.Range("N136").Value = -22938.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1823.2
$ws.Range("I136").Value = 1703.7693
$ws.Range("K136").Value = 5111.3079
$ws.Range("M136").Value = -2561.3079
